# Fill in the "RO & CO Hearing Allocation" worksheet's per-RO allocation
# table (rows 5-60) with the default Video/Virtual scheduling values that
# were previously left at 0 -- mirrors the example values already present
# in the header/example row (row 4): 12 video time slots of 60 minutes
# each, starting at 8:30.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RO & CO Hearing Allocation")

for ($r = 5; $r -le 60; $r++) {
    $ws.Range("F$r").Value = 12
    $ws.Range("G$r").Value = 60
    $ws.Range("H$r").Value = "8:30"
}

# Leave the selection where the user's last edit landed.
$ws.Range("H5").Select()
